# Automation script for test case 2
# - Renames the existing sheet to "ValidLogin" and fills it with the
#   valid-login test data (username/password headers + ADMIN/pointofsale).
# - Adds a second sheet "InvalidLogin" with the invalid-login test data
#   (username/password headers + abcd/xyz) and leaves it as the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

$ws1.Range("A1").Value = "username"
$ws1.Range("B1").Value = "password"
$ws1.Range("B2").Value = "pointofsale"
$ws1.Range("A2").Value = "ADMIN"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

$ws2.Range("B2").Select() | Out-Null
